$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 91-92 (everything from row 91 down shifts down by 2)
$ws.Rows("91:92").Insert()

# Copy the date number format style (style index 2 / $/caja date col) from the row above
$ws.Range("D90").Copy()
$ws.Range("D91:D92").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 91 - new record
$ws.Range("A91").Value = 11
$ws.Range("B91").Value = "Vega Monumental Concepción"
$ws.Range("C91").Value = "Bíobío"
$ws.Range("D91").Value = 45142
$ws.Range("E91").Value = 8
$ws.Range("F91").Value = 100112013
$ws.Range("G91").Value = "Alcachofa"
$ws.Range("H91").Value = "Argentina(o)"
$ws.Range("I91").Value = "Primera"
$ws.Range("J91").Value = 60
$ws.Range("K91").Value = 13000
$ws.Range("L91").Value = 13000
$ws.Range("M91").Value = 13000
$ws.Range("N91").Value = "$/caja 50 unidades"
$ws.Range("O91").Value = "Provincia de Limarí"
$ws.Range("P91").Value = 260
$ws.Range("Q91").Value = 50
$ws.Range("R91").Value = "Hortaliza"

# Row 92 - new record
$ws.Range("A92").Value = 11
$ws.Range("B92").Value = "Vega Monumental Concepción"
$ws.Range("C92").Value = "Bíobío"
$ws.Range("D92").Value = 45142
$ws.Range("E92").Value = 8
$ws.Range("F92").Value = 100112013
$ws.Range("G92").Value = "Alcachofa"
$ws.Range("H92").Value = "Española"
$ws.Range("I92").Value = "Primera"
$ws.Range("J92").Value = 80
$ws.Range("K92").Value = 15000
$ws.Range("L92").Value = 15000
$ws.Range("M92").Value = 15000
$ws.Range("N92").Value = "$/caja 30 unidades"
$ws.Range("O92").Value = "Provincia de Limarí"
$ws.Range("P92").Value = 500
$ws.Range("Q92").Value = 30
$ws.Range("R92").Value = "Hortaliza"
